# Apply updated probability matrix values (games pulled March 7)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1589403973509934
$ws.Range("C2").Value = 0.6456953642384106
$ws.Range("J2").Value = 0.006622516556291391
$ws.Range("P2").Value = 0.1125827814569536
$ws.Range("S2").Value = 0.07615894039735099
$ws.Range("B3").Value = 0.005
$ws.Range("C3").Value = 0.02
$ws.Range("J3").Value = 0.02
$ws.Range("O3").Value = 0.005
$ws.Range("P3").Value = 0.755
$ws.Range("S3").Value = 0.195
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.8292682926829268
$ws.Range("S4").Value = 0.1463414634146341
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.04819277108433735
$ws.Range("D6").Value = 0.01606425702811245
$ws.Range("F6").Value = 0.07630522088353414
$ws.Range("J6").Value = 0.2570281124497992
$ws.Range("O6").Value = 0.0321285140562249
$ws.Range("Q6").Value = 0.1887550200803213
$ws.Range("R6").Value = 0.06827309236947791
$ws.Range("S6").Value = 0.3132530120481928
$ws.Range("B7").Value = 0.1208791208791209
$ws.Range("D7").Value = 0.03846153846153846
$ws.Range("F7").Value = 0.06593406593406594
$ws.Range("J7").Value = 0.1263736263736264
$ws.Range("O7").Value = 0.01098901098901099
$ws.Range("Q7").Value = 0.1978021978021978
$ws.Range("R7").Value = 0.03846153846153846
$ws.Range("S7").Value = 0.4010989010989011
$ws.Range("B8").Value = 0.1047008547008547
$ws.Range("D8").Value = 0.004273504273504274
$ws.Range("F8").Value = 0.07051282051282051
$ws.Range("J8").Value = 0.1538461538461539
$ws.Range("O8").Value = 0.008547008547008548
$ws.Range("Q8").Value = 0.1965811965811966
$ws.Range("R8").Value = 0.1153846153846154
$ws.Range("S8").Value = 0.3461538461538461
$ws.Range("B9").Value = 0.1216931216931217
$ws.Range("D9").Value = 0.01587301587301587
$ws.Range("F9").Value = 0.05291005291005291
$ws.Range("J9").Value = 0.1534391534391534
$ws.Range("O9").Value = 0.01058201058201058
$ws.Range("Q9").Value = 0.1746031746031746
$ws.Range("R9").Value = 0.0582010582010582
$ws.Range("S9").Value = 0.4126984126984127
$ws.Range("B10").Value = 0.1144200626959248
$ws.Range("D10").Value = 0.02037617554858934
$ws.Range("E10").Value = 0.0007836990595611285
$ws.Range("F10").Value = 0.08150470219435736
$ws.Range("J10").Value = 0.1183385579937304
$ws.Range("O10").Value = 0.0164576802507837
$ws.Range("Q10").Value = 0.2210031347962382
$ws.Range("R10").Value = 0.07680250783699059
$ws.Range("S10").Value = 0.3503134796238245
$ws.Range("G11").Value = 0.1328413284132841
$ws.Range("J11").Value = 0.1217712177121771
$ws.Range("K11").Value = 0.1918819188191882
$ws.Range("L11").Value = 0.5535055350553506
$ws.Range("G12").Value = 0.782051282051282
$ws.Range("J12").Value = 0.1794871794871795
$ws.Range("K12").Value = 0.00641025641025641
$ws.Range("L12").Value = 0.02564102564102564
$ws.Range("S12").Value = 0.00641025641025641
$ws.Range("G13").Value = 0.6829268292682927
$ws.Range("J13").Value = 0.2682926829268293
$ws.Range("S13").Value = 0.04878048780487805
$ws.Range("F15").Value = 0.01716738197424893
$ws.Range("H15").Value = 0.2060085836909871
$ws.Range("I15").Value = 0.06008583690987124
$ws.Range("J15").Value = 0.3347639484978541
$ws.Range("K15").Value = 0.06866952789699571
$ws.Range("M15").Value = 0.0128755364806867
$ws.Range("N15").Value = 0.004291845493562232
$ws.Range("O15").Value = 0.0815450643776824
$ws.Range("S15").Value = 0.2145922746781116
$ws.Range("F16").Value = 0.01401869158878505
$ws.Range("H16").Value = 0.191588785046729
$ws.Range("I16").Value = 0.07476635514018691
$ws.Range("J16").Value = 0.3551401869158878
$ws.Range("K16").Value = 0.1121495327102804
$ws.Range("M16").Value = 0.01869158878504673
$ws.Range("O16").Value = 0.0514018691588785
$ws.Range("S16").Value = 0.1822429906542056
$ws.Range("F17").Value = 0.02868852459016394
$ws.Range("H17").Value = 0.1762295081967213
$ws.Range("I17").Value = 0.09836065573770492
$ws.Range("J17").Value = 0.415983606557377
$ws.Range("K17").Value = 0.08811475409836066
$ws.Range("M17").Value = 0.01434426229508197
$ws.Range("O17").Value = 0.06352459016393443
$ws.Range("S17").Value = 0.1147540983606557
$ws.Range("F18").Value = 0.03225806451612903
$ws.Range("H18").Value = 0.2258064516129032
$ws.Range("I18").Value = 0.08064516129032258
$ws.Range("J18").Value = 0.4193548387096774
$ws.Range("K18").Value = 0.06451612903225806
$ws.Range("M18").Value = 0.01075268817204301
$ws.Range("O18").Value = 0.06451612903225806
$ws.Range("S18").Value = 0.1021505376344086
$ws.Range("F19").Value = 0.01827242524916944
$ws.Range("H19").Value = 0.2101328903654485
$ws.Range("I19").Value = 0.08139534883720931
$ws.Range("J19").Value = 0.3612956810631229
$ws.Range("K19").Value = 0.1029900332225914
$ws.Range("M19").Value = 0.02159468438538206
$ws.Range("N19").Value = 0.0008305647840531562
$ws.Range("O19").Value = 0.08139534883720931
$ws.Range("S19").Value = 0.1220930232558139

Write-Output "Applied 110 cell updates to Ohio St._A matrix"
